$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 269
$ws1.Range("F3").Value = 86
$ws1.Range("F4").Value = 964
$ws1.Range("F5").Value = 547

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 269
$ws4.Range("F3").Value = 86
$ws4.Range("F4").Value = 964
$ws4.Range("F6").Value = 547
